$d = $word.ActiveDocument

# Locate the paragraph that currently holds "ANGULAR 2" (and the _GoBack
# bookmark trailing it).
$angularPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "ANGULAR 2*") {
        $angularPara = $p
    }
}

$angularRange = $angularPara.Range
$paraStart = $angularRange.Start
$textEnd = $paraStart + 9            # length of "ANGULAR 2"

# 1) Split the paragraph in two right after "ANGULAR 2": the bookmark
#    (and the paragraph mark it rides on) ends up alone in the new,
#    second paragraph.
$splitPoint = $d.Range($textEnd, $textEnd)
$splitPoint.InsertBefore("`r")

# 2) Replace the "ANGULAR 2" text with two runs: the original text plus a
#    new run holding the parenthetical remark, so it lands in its own
#    <w:r> (as a manually typed follow-up sentence would).
$textRange = $d.Range($paraStart, $textEnd)
$extra = " (MIRAR TRANSPARENCIAS PORQUE HAY UNA MOVIDA QUE NO VEAS MÁQUINA) "
$packageXml = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
  "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
  "<pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
  "<w:body><w:p><w:r><w:t>ANGULAR 2</w:t></w:r><w:r><w:t xml:space=`"preserve`">$extra</w:t></w:r></w:p>" +
  "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$textRange.InsertXML($packageXml)

# 3) Add the new "JAJASALU2" paragraph: its text goes at the very start of
#    the paragraph that now owns the bookmark.
$bookmarkPara = $angularPara.Next()
$bookmarkPara.Range.InsertBefore("JAJASALU2")

Write-Output $d.Content.Text
